# Visualization plan update:
#  - Insert a new block of "planning" bullet paragraphs (Major, Registration,
#    Assignments, Home Pages, Tests, Upload to Azure, If there is time graphs,
#    Presentation training, OOP Principles) at the very top of the document,
#    followed by an empty paragraph that now owns the "_GoBack" bookmark
#    (previously attached to the trailing "Error View" paragraph).
#  - Keep all the pre-existing "TODO:" content below that, unchanged, except
#    for adding a new "Assignments" paragraph right after
#    "The All courses table stays different".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the document, inside
#    the "Error View" paragraph. It needs to move to a new empty paragraph
#    that will be created further down, so remove it from its current spot
#    first.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Insert the new lines before the current first paragraph ("TODO:").
# ---------------------------------------------------------------------------
$newLines = @(
    "Major",
    "Registration",
    "Assignments",
    "Home Pages",
    "Tests",
    "Upload to Azure",
    "If there is time graphs",
    "Presentation training",
    "OOP Principles"
)

$firstPara = $d.Paragraphs(1)
foreach ($line in $newLines) {
    $firstPara.Range.InsertParagraphBefore()
}

for ($i = 1; $i -le $newLines.Count; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.InsertAfter($newLines[$i - 1])
}

# ---------------------------------------------------------------------------
# 3. Insert a brand-new empty paragraph right before "TODO:" and give it the
#    "_GoBack" bookmark (with nothing else inside it).
# ---------------------------------------------------------------------------
$todoPara = $d.Paragraphs($newLines.Count + 1)
$todoPara.Range.InsertParagraphBefore()

$bookmarkPara = $d.Paragraphs($newLines.Count + 1)
# Temporarily add a placeholder character so we can anchor a non-collapsed
# bookmark range fully inside the new paragraph, then remove the character
# again, leaving only the bookmark start/end markers behind.
$bookmarkPara.Range.InsertAfter("X")
$bookmarkStart = $bookmarkPara.Range.Start
$bookmarkRange = $d.Range($bookmarkStart, $bookmarkStart + 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$placeholderRange = $d.Range($bookmarkStart, $bookmarkStart + 1)
$placeholderRange.Text = ""

# ---------------------------------------------------------------------------
# 4. Add the new "Assignments" paragraph right after
#    "The All courses table stays different".
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "The All courses table stays different") {
        $p.Range.InsertParagraphAfter()
        $newAssignmentsPara = $d.Paragraphs($i + 1)
        $newAssignmentsPara.Range.InsertAfter("Assignments")
        break
    }
}
